$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '96.394.42'
$ws.Cells.Item(2, 5).Value = '  -1.22%  '
$ws.Cells.Item(3, 4).Value = '3.328.09'
$ws.Cells.Item(3, 5).Value = '  -2.44%  '
$ws.Cells.Item(4, 5).Value = '  -0.03%  '
$ws.Cells.Item(5, 4).Value = '248.74'
$ws.Cells.Item(5, 5).Value = '  -2.29%  '
$ws.Cells.Item(6, 4).Value = '651.59'
$ws.Cells.Item(6, 5).Value = '  -0.39%  '
$ws.Cells.Item(7, 5).Value = '  -6.38%  '
$ws.Cells.Item(8, 5).Value = '  -1.34%  '
$ws.Cells.Item(9, 4).Value = '1.00'
$ws.Cells.Item(9, 5).Value = '  +0.05%  '
$ws.Cells.Item(10, 4).Value = '0.985'
$ws.Cells.Item(10, 5).Value = '  -7.65%  '
$ws.Cells.Item(11, 4).Value = '3.324.80'
$ws.Cells.Item(11, 5).Value = '  -2.44%  '
$ws.Cells.Item(12, 5).Value = '  -3.49%  '
$ws.Cells.Item(13, 4).Value = '40.22'
$ws.Cells.Item(13, 5).Value = '  -4.17%  '
$ws.Cells.Item(14, 4).Value = '96.091.81'
$ws.Cells.Item(14, 5).Value = '  -1.23%  '
$ws.Cells.Item(15, 5).Value = '  -3.76%  '
$ws.Cells.Item(17, 4).Value = '3.947.43'
$ws.Cells.Item(17, 5).Value = '  -2.52%  '
$ws.Cells.Item(18, 4).Value = '8.49'
$ws.Cells.Item(18, 5).Value = '  -1.93%  '
$ws.Cells.Item(19, 4).Value = '3.318.85'
$ws.Cells.Item(19, 5).Value = '  -2.50%  '
$ws.Cells.Item(20, 5).Value = '  -3.20%  '
$ws.Cells.Item(21, 4).Value = '0.520'
$ws.Cells.Item(21, 5).Value = '  +3.84%  '
$ws.Cells.Item(22, 4).Value = '502.98'
$ws.Cells.Item(22, 5).Value = '  -0.74%  '
$ws.Cells.Item(23, 5).Value = '  -1.86%  '
$ws.Cells.Item(24, 5).Value = '  -3.84%  '
$ws.Cells.Item(25, 5).Value = '  -3.84%  '
$ws.Cells.Item(26, 4).Value = '6.55'
$ws.Cells.Item(26, 5).Value = '  +6.49%  '
$ws.Cells.Item(27, 4).Value = '95.54'
$ws.Cells.Item(27, 5).Value = '  -3.25%  '
$ws.Cells.Item(28, 5).Value = '  -6.09%  '
$ws.Cells.Item(29, 5).Value = '  -8.10%  '
$ws.Cells.Item(30, 5).Value = '  +0.08%  '
$ws.Cells.Item(31, 5).Value = '  -4.10%  '
$ws.Cells.Item(32, 4).Value = '0.187'
$ws.Cells.Item(32, 5).Value = '  -5.89%  '
$ws.Cells.Item(33, 5).Value = '  +8.22%  '
$ws.Cells.Item(34, 4).Value = '1.00'
$ws.Cells.Item(34, 5).Value = '  -0.08%  '
$ws.Cells.Item(35, 4).Value = '0.543'
$ws.Cells.Item(35, 5).Value = '  -5.69%  '
$ws.Cells.Item(36, 5).Value = '  -6.93%  '
$ws.Cells.Item(37, 5).Value = '  +3.28%  '
$ws.Cells.Item(38, 5).Value = '  -1.91%  '
$ws.Cells.Item(39, 5).Value = '  +0.04%  '
$ws.Cells.Item(40, 4).Value = '0.151'
$ws.Cells.Item(40, 5).Value = '  -2.17%  '
$ws.Cells.Item(41, 4).Value = '505.59'
$ws.Cells.Item(41, 5).Value = '  -1.35%  '
$ws.Cells.Item(43, 4).Value = '0.0428'
$ws.Cells.Item(43, 5).Value = '  +0.74%  '
$ws.Cells.Item(44, 5).Value = '  -4.05%  '
$ws.Cells.Item(45, 5).Value = '  -0.54%  '
$ws.Cells.Item(46, 5).Value = '  +5.54%  '
$ws.Cells.Item(47, 4).Value = '5.46'
$ws.Cells.Item(47, 5).Value = '  -1.26%  '
$ws.Cells.Item(48, 5).Value = '  +0.99%  '
$ws.Cells.Item(49, 4).Value = '53.43'
$ws.Cells.Item(49, 5).Value = '  +3.85%  '
$ws.Cells.Item(50, 5).Value = '  -5.04%  '
$ws.Cells.Item(51, 5).Value = '  +1.46%  '
